# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
# Updates the StructureDefinition metadata & root-extension description.

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date refresh
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be the first of two duplicate "Contact" rows;
# it becomes the Jurisdiction row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the second duplicate "Contact" / "No display for ContactDetail"
# row - remove it entirely, shifting everything below up by one.
$meta.Range("A11:B11").EntireRow.Delete()

# --- Sheet "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition now carry the real title & description
# instead of the generic placeholder text.
$elements.Range("K2").Value = "Measure Report Evidence Population Id"
$elements.Range("L2").Value = "Population id of populations that reference or utilize the rule definition"
